$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from the existing
# "sum" header (G1) onto the new "Save" header (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding numeric value (plain, unstyled like F2) in H2.
$ws.Range("H2").Value = 1

$excel.CutCopyMode = $false
